# Localize the Stata "logout" example into Chinese, per the upstream
# diff: heading, command echoes, variable names, and a few regression
# statistics (rounding updates) are all swapped for zh equivalents; two
# intermediate "generate"/"label variable" command-echo lines are
# dropped entirely (the table is described straight off the shipped
# auto_zh dataset, no derived "fuel" variable is created in-document).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        throw "Find.Execute could not find: $find"
    }
}

# 1. Heading
Replace-Text "Include output from Stata commands" "Stata命令输出"

# 2. The "sysuse" command echo now loads the Chinese-labelled dataset
Replace-Text ". sysuse auto, clear" ". sysuse auto_zh, clear"

# 3. Dataset description echo
Replace-Text "(1978 Automobile Data)" "(1978年汽车数据)"

# 4. Drop the "generate fuel ..." / "label variable fuel ..." echoed
#    lines (and their blank-line separators) entirely -- `v is a
#    vertical-tab, the character Word's Range.Text uses for a
#    manual line break (<w:br/>).
$genLabelBlock = ". generate fuel = 100/mpg`v`v. label variable fuel `"Fuel consumption (Gallons per 100 Miles)`"`v`v"
Replace-Text $genLabelBlock ""

# 5. "describe" now targets the renamed (already-Chinese) variables
Replace-Text ". describe fuel weight" ". describe 油耗 重量"

# 6. The two-line "fuel" describe row (wrapped across lines) collapses
#    into a single "油耗" line, and the "weight" row becomes "重量"
#    (also single line each now).
$describeBlock = "fuel            float   %9.0g                 Fuel consumption (Gallons per 100`v                                                Miles)`vweight          int     %8.0gc                Weight (lbs.)"
$describeReplacement = "油耗            float   %9.0g                 油量消耗(公升每一百公里)`v重量            float   %8.0gc                重量(公斤)"
Replace-Text $describeBlock $describeReplacement

# 7. "regress" command echo
Replace-Text ". regress fuel weight" ". regress 油耗 重量"

# 8. ANOVA "Model" row -- updated rounding (...4969 -> ...4971)
Replace-Text "       Model |  87.2964969         1  87.2964969   Prob > F        =    0.0000" "       Model |  87.2964971         1  87.2964971   Prob > F        =    0.0000"

# 9. ANOVA "Residual" row -- updated rounding (...639 -> ...637 / ...054 -> ...051)
Replace-Text "    Residual |  32.2797639        72  .448330054   R-squared       =    0.7300" "    Residual |  32.2797637        72  .448330051   R-squared       =    0.7300"

# 10. Coefficient table header ("fuel" -> "油耗")
Replace-Text "        fuel |      Coef.   Std. Err.      t    P>|t|     [95% Conf. Interval]" "        油耗 |      Coef.   Std. Err.      t    P>|t|     [95% Conf. Interval]"

# 11. Coefficient table "weight" row -> "重量" row with rescaled coefficients
Replace-Text "      weight |    .001407   .0001008    13.95   0.000      .001206    .0016081" "        重量 |    .003102   .0002223    13.95   0.000     .0026589    .0035452"

Write-Output "done"
